# Ocean accounts to RDF/TTL
#
# Updates the SSSOM sheet's "O" column formulas so the generated label also
# includes the row number of the mapping (via ROW()-1), and leaves the
# workbook with the SSSOM sheet active/selected (mirroring the author's
# final on-screen state when the file was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SSSOM")

# Extend each "<Ecosystem> - mapping to IUCN GET" formula with the mapping's
# row number, e.g. "Saltmarsh - mapping to IUCN GET - 1"
$ws.Range("O2").Formula = '=_xlfn.CONCAT(B2, " - mapping to IUCN GET - ", ROW(B2)-1 )'
$ws.Range("O3").Formula = '=_xlfn.CONCAT(B3, " - mapping to IUCN GET - ", ROW(B3)-1 )'
$ws.Range("O4").Formula = '=_xlfn.CONCAT(B4, " - mapping to IUCN GET - ", ROW(B4)-1 )'
$ws.Range("O5").Formula = '=_xlfn.CONCAT(B5, " - mapping to IUCN GET - ", ROW(B5)-1 )'

# Make SSSOM the active/selected sheet (it becomes the tab shown when the
# workbook is reopened), with O3:O5 selected and O3 as the active cell.
$ws.Activate()
$ws.Range("O3:O5").Select()
